$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = -0.0003960461763199419
$ws.Range("F2").Value = 71.58631896972656
$ws.Range("G2").Value = [double]"-4.48198989033699e-09"
$ws.Range("H2").Value = [double]"2.980232238769531e-08"
$ws.Range("I2").Value = 0.0863189697265625
$ws.Range("J2").Value = [double]"-4.48198989033699e-09"
$ws.Range("K2").Value = 0.0003960759786423296
$ws.Range("L2").Value = 0.08631987842227357
$ws.Range("E3").Value = -0.0003960461763199419
$ws.Range("F3").Value = -71.58631896972656
$ws.Range("G3").Value = [double]"4.423782229423523e-09"
$ws.Range("I3").Value = -0.0863189697265625
$ws.Range("J3").Value = [double]"4.423782229423523e-09"
$ws.Range("K3").Value = 0.0003960461763199419
$ws.Range("L3").Value = 0.08631987828553173
$ws.Range("E4").Value = 118.8882751464844
$ws.Range("F4").Value = 318.4662780761719
$ws.Range("G4").Value = 390.2781372070312
$ws.Range("H4").Value = 115.4920272827148
$ws.Range("I4").Value = -3.033721923828125
$ws.Range("J4").Value = 4.42877197265625
$ws.Range("K4").Value = -3.396247863769531
$ws.Range("L4").Value = 6.352321579474709
$ws.Range("E5").Value = 118.8882751464844
$ws.Range("F5").Value = -327.5456848144531
$ws.Range("G5").Value = 390.4197387695312
$ws.Range("H5").Value = 116.120231628418
$ws.Range("I5").Value = -6.045684814453125
$ws.Range("J5").Value = 4.57037353515625
$ws.Range("K5").Value = -2.768043518066406
$ws.Range("L5").Value = 8.068499491508641
$ws.Range("E6").Value = 118.8882751464844
$ws.Range("F6").Value = 327.8734741210938
$ws.Range("G6").Value = -389.3250122070312
$ws.Range("H6").Value = 115.4021759033203
$ws.Range("I6").Value = 6.37347412109375
$ws.Range("J6").Value = -3.47564697265625
$ws.Range("K6").Value = -3.486099243164062
$ws.Range("L6").Value = 8.053209433758404
$ws.Range("E7").Value = 118.8882751464844
$ws.Range("F7").Value = -319.3450012207031
$ws.Range("G7").Value = -392.5592346191406
$ws.Range("H7").Value = 115.7499008178711
$ws.Range("I7").Value = 2.154998779296875
$ws.Range("J7").Value = -6.709869384765625
$ws.Range("K7").Value = -3.138374328613281
$ws.Range("L7").Value = 7.714645832822458
$ws.Range("E8").Value = 125.0787048339844
$ws.Range("G8").Value = 1.454944252967834
$ws.Range("H8").Value = 121.307746887207
$ws.Range("J8").Value = 1.454944252967834
$ws.Range("K8").Value = -3.770957946777344
$ws.Range("L8").Value = 5.205399062862468
$ws.Range("E9").Value = 125.0787048339844
$ws.Range("G9").Value = -1.736018896102905
$ws.Range("H9").Value = 122.4498672485352
$ws.Range("J9").Value = -1.736018896102905
$ws.Range("K9").Value = -2.628837585449219
$ws.Range("L9").Value = 4.654996845224717
$ws.Range("E10").Value = -6.126372814178467
$ws.Range("F10").Value = -429.8712244205686
$ws.Range("G10").Value = 383.3324926272199
$ws.Range("H10").Value = -7.909320997443928
$ws.Range("I10").Value = -6.602655326716615
$ws.Range("J10").Value = 1.495883740501142
$ws.Range("K10").Value = -1.782948183265461
$ws.Range("L10").Value = 7.000830647339257
$ws.Range("E11").Value = -6.126372814178467
$ws.Range("F11").Value = -503.2148461615025
$ws.Range("G11").Value = 383.3324926272199
$ws.Range("H11").Value = -7.909320997443928
$ws.Range("I11").Value = -8.446277067650499
$ws.Range("J11").Value = 1.495883740501142
$ws.Range("K11").Value = -1.782948183265461
$ws.Range("L11").Value = 8.761059792788997
$ws.Range("E12").Value = 125.0787048339844
$ws.Range("F12").Value = -503.2148461615025
$ws.Range("G12").Value = -2.083690761659004
$ws.Range("H12").Value = 120.4774581693322
$ws.Range("I12").Value = -8.446291900221013
$ws.Range("J12").Value = -2.083683132264472
$ws.Range("K12").Value = -4.60124666465218
$ws.Range("L12").Value = 9.84139487717038
$ws.Range("E13").Value = -6.126372814178467
$ws.Range("F13").Value = -425.284281945186
$ws.Range("G13").Value = -390.438538203389
$ws.Range("H13").Value = -5.937114971814367
$ws.Range("I13").Value = -2.015712851333944
$ws.Range("J13").Value = -8.601929316670294
$ws.Range("K13").Value = 0.1892578423640998
$ws.Range("L13").Value = 8.836973735330496
$ws.Range("E14").Value = -6.126372814178467
$ws.Range("F14").Value = -503.2148461615025
$ws.Range("G14").Value = -390.438538203389
$ws.Range("H14").Value = -5.937114971814367
$ws.Range("I14").Value = -8.446277067650499
$ws.Range("J14").Value = -8.601929316670294
$ws.Range("K14").Value = 0.1892578423640998
$ws.Range("L14").Value = 12.0568902625597
$ws.Range("E15").Value = -6.126372814178467
$ws.Range("F15").Value = 502.6868610586105
$ws.Range("G15").Value = 391.1159981042497
$ws.Range("H15").Value = -10.50787940809246
$ws.Range("I15").Value = 7.918291964758509
$ws.Range("J15").Value = 9.279389217530934
$ws.Range("K15").Value = -4.381506593913997
$ws.Range("L15").Value = 12.96163615914674
$ws.Range("E16").Value = -6.126372814178467
$ws.Range("F16").Value = 421.8202224262885
$ws.Range("G16").Value = 391.1159981042497
$ws.Range("H16").Value = -10.50787940809246
$ws.Range("I16").Value = -1.448346667563499
$ws.Range("J16").Value = 9.279389217530934
$ws.Range("K16").Value = -4.381506593913997
$ws.Range("L16").Value = 10.36351158403287
$ws.Range("E17").Value = 125.0787048339844
$ws.Range("F17").Value = 502.6868610586105
$ws.Range("G17").Value = 9.25552607000553
$ws.Range("H17").Value = 119.0693511043597
$ws.Range("I17").Value = 7.918306797329024
$ws.Range("J17").Value = 9.255533699400061
$ws.Range("K17").Value = -6.009353729624678
$ws.Range("L17").Value = 13.58222437029779
$ws.Range("E18").Value = -6.126372814178467
$ws.Range("F18").Value = 502.6868610586105
$ws.Range("G18").Value = -381.0895568217501
$ws.Range("H18").Value = -3.234903379534614
$ws.Range("I18").Value = 7.918291964758509
$ws.Range("J18").Value = 0.7470520649686705
$ws.Range("K18").Value = 2.891469434643852
$ws.Range("L18").Value = 8.462743640121257
$ws.Range("E19").Value = -6.126372814178467
$ws.Range("F19").Value = 436.3910100070559
$ws.Range("G19").Value = -381.0895568217501
$ws.Range("H19").Value = -3.234903379534614
$ws.Range("I19").Value = 13.12244091320383
$ws.Range("J19").Value = 0.7470520649686705
$ws.Range("K19").Value = 2.891469434643852
$ws.Range("L19").Value = 13.45797673499919
